$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C (Förändrad) rows 2-264 all change from serial date 45188 to 45189
$ws.Range("C2:C264").Value = 45189
